$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Personal Website" column header to "Google Scholar Profile"
$ws.Range("D1").Value = "Google Scholar Profile"

# Remove the last paper row (row 8: "Multi-agent Architecture Search via Agentic Supernet")
$ws.Rows("8:8").Delete()

# Recompute the sheet's outline level metadata after the row removal
$ws.Outline.ShowLevels(6, 3) | Out-Null

# Update the saved selection/active cell
$ws.Range("C13").Select() | Out-Null
